$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprints were previously labelled as 1 week each (Sprintti 1..5) but are
# actually 2 weeks long. Relabel the sprint headers so two consecutive
# date ranges share the same sprint number.
$ws.Range("A14").Value = "Sprintti 2"
$ws.Range("A20").Value = "Sprintti 3"
$ws.Range("A26").Value = "Sprintti 3"

# Update the active selection left behind by the author.
$ws.Range("B11:G11").Select()
